$d = $word.ActiveDocument

# Replace the literal storage-section intro sentence with the [storageintro] placeholder.
$d.Content.Find.Execute(
    "For the duration of the project, storage and backup of data will be ensured by the project manager.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "[storageintro]", 2
)
